# Refresh the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E) to match the latest scrape. Values are written
# as literal text (matching the sheet's existing inlineStr cells) --
# numeric-looking prices are prefixed with a leading apostrophe so Excel
# keeps them as text instead of coercing them to floating-point numbers
# (which would corrupt the fixed-decimal formatting, e.g. "213.33" ->
# 213.3300000000000).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.260.82'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.592.09'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''213.33'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('D10').Value = '''18.97'
$ws.Range('E10').Value = '  -2.23%  '
$ws.Range('D11').Value = '''0.0851'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '1.816.08'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = '1.595.90'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('E14').Value = '  -1.19%  '
$ws.Range('E15').Value = '  -2.61%  '
$ws.Range('D16').Value = '''63.79'
$ws.Range('E16').Value = '  -1.26%  '
$ws.Range('D17').Value = '26.245.03'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').Value = '''215.12'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('E20').Value = '  -1.68%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').Value = '''4.28'
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('D23').Value = '''9.04'
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('E24').Value = '  -3.22%  '
$ws.Range('D25').Value = '''144.76'
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '''6.96'
$ws.Range('E27').Value = '  -1.47%  '
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('D29').Value = '''15.10'
$ws.Range('E29').Value = '  -0.78%  '
$ws.Range('D30').Value = '''0.0490'
$ws.Range('E30').Value = '  -2.36%  '
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('E32').Value = '  -0.64%  '
$ws.Range('D33').Value = '1.409.25'
$ws.Range('E33').Value = '  +5.19%  '
$ws.Range('E34').Value = '  -0.14%  '
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('E36').Value = '  -1.58%  '
$ws.Range('D37').Value = '''0.575'
$ws.Range('E37').Value = '  -4.29%  '
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('E39').Value = '  +0.43%  '
$ws.Range('E40').Value = '  -0.66%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').Value = '''0.959'
$ws.Range('E42').Value = '  -7.71%  '
$ws.Range('E43').Value = '  +1.11%  '
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').Value = '1.728.33'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = '''60.85'
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').Value = '''87.05'
$ws.Range('E47').Value = '  -0.83%  '
$ws.Range('E48').Value = '  -1.35%  '
$ws.Range('D49').Value = '''0.0502'
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('E50').Value = '  -3.45%  '
$ws.Range('E51').Value = '  +0.06%  '
